$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$headers = @("id","name","brand","model","mtyp_code","min_driver_ver","descr","lang_code","is_active","cr_by","cr_dtimes","upd_by","upd_dtimes","is_deleted","del_dtimes")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data row (row 2)
$ws.Cells.Item(2, 1).Value = 1001
$ws.Cells.Item(2, 2).Value = "PROBOOK"
$ws.Cells.Item(2, 3).Value = "HP"
$ws.Cells.Item(2, 4).Value = "430 G6"
$ws.Cells.Item(2, 5).Value = "LTP"
$ws.Cells.Item(2, 6).Value = 1.454
$ws.Cells.Item(2, 7).Value = "PC pour pour les enrÃ´lements en phase developpement"
$ws.Cells.Item(2, 8).Value = "fra"
$ws.Cells.Item(2, 9).Value = $true
$ws.Cells.Item(2, 10).Value = "superadmin"
$ws.Cells.Item(2, 11).NumberFormat = "mm:ss.0"
$ws.Cells.Item(2, 11).Value = 45079.57824934028
$ws.Cells.Item(2, 12).Value = "NULL"
$ws.Cells.Item(2, 13).Value = "NULL"
$ws.Cells.Item(2, 14).Value = $false
$ws.Cells.Item(2, 15).Value = "NULL"

$ws.Range("A1:O2").Select()
